$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4062.625
$ws.Range("I40").Value = 2833.6667
$ws.Range("K40").Value = 2833.6667
$ws.Range("M40").Value = -2658.6667

$ws.Range("H70").Value = 1422
$ws.Range("I70").Value = 1466.6666
$ws.Range("J70").Value = 1399.6666
$ws.Range("K70").Value = 4399.9998
$ws.Range("L70").Value = 4198.9998
$ws.Range("M70").Value = -4129.9998
$ws.Range("N70").Value = -4738.9998

$ws.Range("H73").Value = 1422
$ws.Range("I73").Value = 1466.6666
$ws.Range("J73").Value = 1399.6666
$ws.Range("K73").Value = 4399.9998
$ws.Range("L73").Value = 4198.9998
$ws.Range("M73").Value = -3463.9998
$ws.Range("N73").Value = -6070.9998

$ws.Range("H74").Value = 4042.8
$ws.Range("I74").Value = 3182.4443
$ws.Range("J74").Value = 5333.3335
$ws.Range("K74").Value = 3182.4443
$ws.Range("L74").Value = 5333.3335
$ws.Range("M74").Value = -2246.4443
$ws.Range("N74").Value = -7205.3335

$ws.Range("H76").Value = 3877.75
$ws.Range("I76").Value = 3307.25
$ws.Range("J76").Value = 5018.75
$ws.Range("K76").Value = 3307.25
$ws.Range("L76").Value = 5018.75
$ws.Range("M76").Value = -2992.25
$ws.Range("N76").Value = -5648.75

$ws.Range("H77").Value = 4042.8
$ws.Range("I77").Value = 3182.4443
$ws.Range("J77").Value = 5333.3335
$ws.Range("K77").Value = 15912.2215
$ws.Range("L77").Value = 26666.6675
$ws.Range("M77").Value = -11232.2215
$ws.Range("N77").Value = -36026.6675

$ws.Range("H79").Value = 3877.75
$ws.Range("I79").Value = 3307.25
$ws.Range("J79").Value = 5018.75
$ws.Range("K79").Value = 3307.25
$ws.Range("L79").Value = 5018.75
$ws.Range("M79").Value = -2215.25
$ws.Range("N79").Value = -7202.75

$ws.Range("H86").Value = 3049.75
$ws.Range("I86").Value = 3099.5
$ws.Range("K86").Value = 3099.5
$ws.Range("M86").Value = -1976.5

$ws.Range("H89").Value = 3049.75
$ws.Range("I89").Value = 3099.5
$ws.Range("K89").Value = 15497.5
$ws.Range("M89").Value = -9881.5

$ws.Range("H113").Value = 5512.1875
$ws.Range("J113").Value = 5039.6
$ws.Range("L113").Value = 5039.6
$ws.Range("N113").Value = -11547.6

$ws.Range("H116").Value = 4639.154
$ws.Range("I116").Value = 4761.8
$ws.Range("J116").Value = 4562.5
$ws.Range("K116").Value = 4761.8
$ws.Range("L116").Value = 4562.5
$ws.Range("M116").Value = -1319.8
$ws.Range("N116").Value = -11446.5

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H129").Value = 1716.7333
$ws.Range("I129").Value = 1670.8182
$ws.Range("J129").Value = 1843
$ws.Range("K129").Value = 5012.4546
$ws.Range("L129").Value = 5529
$ws.Range("M129").Value = -12.45460000000003
$ws.Range("N129").Value = -15529

$ws.Range("H137").Value = 4785.25
$ws.Range("I137").Value = 1649.3793
$ws.Range("K137").Value = 4948.1379
$ws.Range("M137").Value = -2398.1379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36790.484
$ws.Range("I32").Value = 20892.77
$ws.Range("J32").Value = 140125.62
$ws.Range("K32").Value = 20892.77
$ws.Range("L32").Value = 140125.62
$ws.Range("M32").Value = -20605.77
$ws.Range("N32").Value = -140699.62

$ws.Range("H122").Value = 1523.8572
$ws.Range("I122").Value = 1379.3793
$ws.Range("J122").Value = 2222.1667
$ws.Range("K122").Value = 4138.1379
$ws.Range("L122").Value = 6666.500100000001
$ws.Range("M122").Value = -1688.1379
$ws.Range("N122").Value = -11566.5001

$ws.Range("H132").Value = 2189.16
$ws.Range("I132").Value = 1711.7646
$ws.Range("J132").Value = 3203.625
$ws.Range("K132").Value = 5135.293799999999
$ws.Range("L132").Value = 9610.875
$ws.Range("M132").Value = -2605.293799999999
$ws.Range("N132").Value = -14670.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3824.6428
$ws.Range("I137").Value = 1657.875
$ws.Range("K137").Value = 4973.625
$ws.Range("M137").Value = 126.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3515.5833
$ws.Range("I132").Value = 3470.2856
$ws.Range("J132").Value = 3579
$ws.Range("K132").Value = 10410.8568
$ws.Range("L132").Value = 10737
$ws.Range("M132").Value = -7880.856800000001
$ws.Range("N132").Value = -15797

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19816
$ws.Range("I7").Value = 38341.285
$ws.Range("J7").Value = 3606.375
$ws.Range("K7").Value = 38341.285
$ws.Range("L7").Value = 3606.375
$ws.Range("M7").Value = -38229.285
$ws.Range("N7").Value = -3830.375

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H46").Value = 1799.5
$ws.Range("I46").Value = 1719.5
$ws.Range("K46").Value = 1719.5
$ws.Range("M46").Value = -1531.5

$ws.Range("H87").Value = 25189
$ws.Range("J87").Value = 25189
$ws.Range("L87").Value = 25189
$ws.Range("N87").Value = -27435

$ws.Range("H90").Value = 25189
$ws.Range("J90").Value = 25189
$ws.Range("L90").Value = 75567
$ws.Range("N90").Value = -86799

$ws.Range("H93").Value = 2279.7222
$ws.Range("I93").Value = 2181.6365
$ws.Range("J93").Value = 2433.8572
$ws.Range("K93").Value = 2181.6365
$ws.Range("L93").Value = 2433.8572
$ws.Range("M93").Value = -933.6365000000001
$ws.Range("N93").Value = -4929.8572

$ws.Range("H126").Value = 19816
$ws.Range("I126").Value = 38341.285
$ws.Range("J126").Value = 3606.375
$ws.Range("K126").Value = 115023.855
$ws.Range("L126").Value = 10819.125
$ws.Range("M126").Value = -112553.855
$ws.Range("N126").Value = -15759.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5555.5557
$ws.Range("I81").Value = 5555.5557
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 11111.1114
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -10050.1114
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 5555.5557
$ws.Range("I84").Value = 5555.5557
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 55555.557
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -50251.557
$ws.Range("N84").ClearContents()

$ws.Range("H107").Value = 62501576
$ws.Range("I107").Value = 1683.3334
$ws.Range("J107").Value = 250001250
$ws.Range("K107").Value = 5050.0002
$ws.Range("L107").Value = 750003750
$ws.Range("M107").Value = -3130.0002
$ws.Range("N107").Value = -750007590
